$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I20").Value = -1.090592475347884
$ws.Range("J20").Value = 0.305552710783839
$ws.Range("K20").Value = 0.459843913720376
$ws.Range("L20").Value = 2.595473261476525
